# Update cohort and classifier QA table (Ascertainment_Overlap_Intersections)
# - Refreshes the ICD/ABG/VBG/OTHER boolean flag combinations and counts in
#   rows 2-13
# - Appends three new combination rows (14-16)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = $False
$ws.Cells.Item(2, 2).Value = $True
$ws.Cells.Item(2, 3).Value = $False
$ws.Cells.Item(2, 4).Value = $False
$ws.Cells.Item(2, 5).Value = 4333
$ws.Cells.Item(3, 1).Value = $False
$ws.Cells.Item(3, 2).Value = $False
$ws.Cells.Item(3, 3).Value = $True
$ws.Cells.Item(3, 4).Value = $False
$ws.Cells.Item(3, 5).Value = 3041
$ws.Cells.Item(4, 1).Value = $False
$ws.Cells.Item(4, 2).Value = $True
$ws.Cells.Item(4, 3).Value = $True
$ws.Cells.Item(4, 4).Value = $False
$ws.Cells.Item(4, 5).Value = 1392
$ws.Cells.Item(5, 1).Value = $True
$ws.Cells.Item(5, 2).Value = $True
$ws.Cells.Item(5, 3).Value = $True
$ws.Cells.Item(5, 4).Value = $False
$ws.Cells.Item(5, 5).Value = 536
$ws.Cells.Item(6, 1).Value = $True
$ws.Cells.Item(6, 2).Value = $False
$ws.Cells.Item(6, 3).Value = $True
$ws.Cells.Item(6, 4).Value = $False
$ws.Cells.Item(6, 5).Value = 467
$ws.Cells.Item(7, 1).Value = $True
$ws.Cells.Item(7, 2).Value = $False
$ws.Cells.Item(7, 3).Value = $False
$ws.Cells.Item(7, 4).Value = $False
$ws.Cells.Item(7, 5).Value = 460
$ws.Cells.Item(8, 1).Value = $False
$ws.Cells.Item(8, 2).Value = $True
$ws.Cells.Item(8, 3).Value = $True
$ws.Cells.Item(8, 4).Value = $True
$ws.Cells.Item(8, 5).Value = 327
$ws.Cells.Item(9, 1).Value = $False
$ws.Cells.Item(9, 2).Value = $False
$ws.Cells.Item(9, 3).Value = $False
$ws.Cells.Item(9, 4).Value = $True
$ws.Cells.Item(9, 5).Value = 262
$ws.Cells.Item(10, 1).Value = $False
$ws.Cells.Item(10, 2).Value = $True
$ws.Cells.Item(10, 3).Value = $False
$ws.Cells.Item(10, 4).Value = $True
$ws.Cells.Item(10, 5).Value = 247
$ws.Cells.Item(11, 1).Value = $True
$ws.Cells.Item(11, 2).Value = $True
$ws.Cells.Item(11, 3).Value = $True
$ws.Cells.Item(11, 4).Value = $True
$ws.Cells.Item(11, 5).Value = 226
$ws.Cells.Item(12, 1).Value = $True
$ws.Cells.Item(12, 2).Value = $True
$ws.Cells.Item(12, 3).Value = $False
$ws.Cells.Item(12, 4).Value = $False
$ws.Cells.Item(12, 5).Value = 190
$ws.Cells.Item(13, 1).Value = $False
$ws.Cells.Item(13, 2).Value = $False
$ws.Cells.Item(13, 3).Value = $True
$ws.Cells.Item(13, 4).Value = $True
$ws.Cells.Item(13, 5).Value = 184
$ws.Cells.Item(14, 1).Value = $True
$ws.Cells.Item(14, 2).Value = $False
$ws.Cells.Item(14, 3).Value = $True
$ws.Cells.Item(14, 4).Value = $True
$ws.Cells.Item(14, 5).Value = 71
$ws.Cells.Item(15, 1).Value = $True
$ws.Cells.Item(15, 2).Value = $True
$ws.Cells.Item(15, 3).Value = $False
$ws.Cells.Item(15, 4).Value = $True
$ws.Cells.Item(15, 5).Value = 19
$ws.Cells.Item(16, 1).Value = $True
$ws.Cells.Item(16, 2).Value = $False
$ws.Cells.Item(16, 3).Value = $False
$ws.Cells.Item(16, 4).Value = $True
$ws.Cells.Item(16, 5).Value = 14
